# The data table used to start at row 1. To leave room for a couple of
# header/title rows (matching the refactored Read-ExcelData function's new
# "dynamic starting row" support), insert two blank rows above the table so
# the existing data shifts down to rows 3:28.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("1:2").Insert()

# Leave the selection on the newly inserted (still blank) second row, the
# way Excel leaves the selection after an insert-rows operation performed
# with the row header selected.
$ws.Rows("2:2").Select()
